$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions) - F column "想去人数" (want-to-go count) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 202
$ws1.Range("F3").Value = 180
$ws1.Range("F4").Value = 1348
$ws1.Range("F5").Value = 18809
$ws1.Range("F6").Value = 403
$ws1.Range("F7").Value = 281
$ws1.Range("F9").Value = 6999
$ws1.Range("F10").Value = 440
$ws1.Range("F11").Value = 701
$ws1.Range("F12").Value = 174
$ws1.Range("F13").Value = 20
$ws1.Range("F14").Value = 130
$ws1.Range("F15").Value = 81
$ws1.Range("F16").Value = 226
$ws1.Range("F17").Value = 175
$ws1.Range("F19").Value = 303
$ws1.Range("F21").Value = 668
$ws1.Range("F23").Value = 40
$ws1.Range("F24").Value = 42
$ws1.Range("F25").Value = 289
$ws1.Range("F26").Value = 1031
$ws1.Range("F27").Value = 9
$ws1.Range("F30").Value = 544
$ws1.Range("F31").Value = 22
$ws1.Range("F32").Value = 88
$ws1.Range("F35").Value = 12265
$ws1.Range("F36").Value = 1301
$ws1.Range("F37").Value = 29
$ws1.Range("F38").Value = 47
$ws1.Range("F39").Value = 225
$ws1.Range("F40").Value = 310
$ws1.Range("F41").Value = 3946
$ws1.Range("F42").Value = 307

# Sheet 2: 演出 (Performances) - F column updates
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 2

# Sheet 4: 全部类型 (All types) - F column updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 180
$ws4.Range("F4").Value = 1348
$ws4.Range("F5").Value = 18809
$ws4.Range("F6").Value = 403
$ws4.Range("F7").Value = 281
$ws4.Range("F9").Value = 6999
$ws4.Range("F10").Value = 440
$ws4.Range("F11").Value = 701
$ws4.Range("F12").Value = 174
$ws4.Range("F13").Value = 20
$ws4.Range("F14").Value = 130
$ws4.Range("F15").Value = 81
$ws4.Range("F16").Value = 226
$ws4.Range("F17").Value = 175
$ws4.Range("F18").Value = 1312
$ws4.Range("F19").Value = 303
$ws4.Range("F21").Value = 668
$ws4.Range("F23").Value = 40
$ws4.Range("F24").Value = 42
$ws4.Range("F26").Value = 1031
$ws4.Range("F27").Value = 9
$ws4.Range("F29").Value = 5204
$ws4.Range("F30").Value = 544
$ws4.Range("F31").Value = 2
$ws4.Range("F32").Value = 22
$ws4.Range("F34").Value = 88
$ws4.Range("F35").Value = 18
$ws4.Range("F36").Value = 80
$ws4.Range("F37").Value = 12265
$ws4.Range("F38").Value = 1301
$ws4.Range("F39").Value = 29
$ws4.Range("F40").Value = 47
$ws4.Range("F41").Value = 225
$ws4.Range("F42").Value = 310
$ws4.Range("F44").Value = 307
